# Generate Report for Handoff
# Update the "latest" datetime values for the last row (d237c79e-... entry)
# on the Overview, zh-cn and de-de sheets to reflect the new handoff/handback.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: column D = "Latest Handoff Date" for row 7 (d237c79e-... file)
$overview.Range("D7").Value = "2016-48-13 22:48:53"

# zh-cn sheet: column E = "Latest Handoff Datetime" for row 7 (d237c79e-... file)
$zhcn.Range("E7").Value = "2016-03-13 22:48:49"

# de-de sheet: column E = "Latest Handoff Datetime" for row 7 (d237c79e-... file)
$dede.Range("E7").Value = "2016-03-13 22:48:53"
